$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.873.61"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").Value = "1.563.49"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "205.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.487"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -1.34%  "
$ws.Range("E9").Value = "  -0.28%  "
$ws.Range("E10").Value = "  -0.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0865"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("D12").Value = "1.785.91"
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("D13").Value = "1.563.51"
$ws.Range("E13").Value = "  +0.36%  "
$ws.Range("E14").Value = "  -0.80%  "
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").Value = "26.874.83"
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.28"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "215.32"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.31%  "
$ws.Range("D20").Value = "0.0₃0683"
$ws.Range("E20").Value = "  -0.40%  "
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.75%  "
$ws.Range("E26").Value = "  +1.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("E29").Value = "  -0.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0465"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.13%  "
$ws.Range("E31").Value = "  -3.41%  "
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("D33").Value = "1.392.62"
$ws.Range("E33").Value = "  +0.71%  "
$ws.Range("E34").Value = "  -0.33%  "
$ws.Range("E35").Value = "  -1.12%  "
$ws.Range("E36").Value = "  -0.54%  "
$ws.Range("E37").Value = "  -2.35%  "
$ws.Range("E38").Value = "  -0.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.530"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.812"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.40%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("E42").Value = "  +0.37%  "
$ws.Range("E43").Value = "  +5.72%  "
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("E45").Value = "  +1.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.70"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.61%  "
$ws.Range("D47").Value = "1.699.86"
$ws.Range("E47").Value = "  +0.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.78"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.66%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₇0983"
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0503"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.40%  "
$ws.Range("E51").Value = "  +1.19%  "
